$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Normalize formatting on the existing rows (2-12).
#    Column A -> default "Normal" style (no number format, no hyperlink font)
#    Column B -> currency style (numFmt 6), except row 9 which keeps the
#                2-decimal currency style (numFmt 8)
#    Column C -> "Hyperlink" named style (keeps the existing hyperlink link)
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Style = "Normal"
    if ($r -eq 9) {
        $ws.Cells.Item($r, 2).NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
    } else {
        $ws.Cells.Item($r, 2).NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'
    }
    $ws.Cells.Item($r, 3).Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 2. Add the three new BOM rows for the tread kit.
# ---------------------------------------------------------------------------

# Row 13 - Tamiya Track and Wheel Set
$ws.Range("A13").Value = "Tamiya Track and Wheel Set"
$ws.Range("A13").Style = "Normal"
$ws.Range("B13").Value = 13
$ws.Range("B13").NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'
$ws.Range("C13").Value = "https://www.amazon.com/Tamiya-Track-Wheel-Set-70100/dp/B001VZJDY2/ref=sr_1_1?keywords=Tamiya+Track+and+Wheel+Set&sr=8-1 "
$ws.Hyperlinks.Add($ws.Range("C13"), "https://www.amazon.com/Tamiya-Track-Wheel-Set-70100/dp/B001VZJDY2/ref=sr_1_1?keywords=Tamiya+Track+and+Wheel+Set&sr=8-1")

# Row 14 - Breadboard
$ws.Range("A14").Value = "Breadboard"
$ws.Range("A14").Style = "Normal"
$ws.Range("B14").Value = 2
$ws.Range("B14").NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'
$breadboardUrl = "https://www.amazon.com/Pcs-MCIGICM-Points-Solderless-Breadboard/dp/B07PCJP9DY/ref=sr_1_1_sspa?keywords=400+Tie-points+breadboard&sr=8-1-spons&psc=1&spLa=ZW5jcnlwdGVkUXVhbGlmaWVyPUEyTVRVVTE1NjBHVzVQJmVuY3J5cHRlZElkPUEwMjkzNDIxMTNHTTNMWDJENEw3NCZlbmNyeXB0ZWRBZElkPUEwMDM1MDQ2Mk4zR1VZSVU0NTJBMiZ3aWRnZXROYW1lPXNwX2F0ZiZhY3Rpb249Y2xpY2tSZWRpcmVjdCZkb05vdExvZ0NsaWNrPXRydWU="
$ws.Hyperlinks.Add($ws.Range("C14"), $breadboardUrl, "", "", ($breadboardUrl + " "))

# Row 15 - Basswood
$ws.Range("A15").Value = "Basswood"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = 15
$ws.Range("B15").NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'
$basswoodUrl = "https://www.amazon.com/Hammont-Basswood-Sheets-12x8x1-Pack/dp/B09NWBMDNL/ref=sxin_17_ac_d_mf_brs?ac_md=3-1-SGFtbW9udA%3D%3D-ac_d_mf_brs_brs&content-id=amzn1.sym.1ad31f34-ba12-4dca-be4b-f62f7f5bb10d%3Aamzn1.sym.1ad31f34-ba12-4dca-be4b-f62f7f5bb10d&cv_ct_cx=basswood&keywords=basswood&pd_rd_i=B09NWBMDNL&pd_rd_r=5220cbe5-6806-437d-863d-ca70e145ba67&pd_rd_w=LW4Pl&pd_rd_wg=ExU61&pf_rd_p=1ad31f34-ba12-4dca-be4b-f62f7f5bb10d&pf_rd_r=70TFQK3M4HVSSSABBFN7&sbo=RZvfv%2F%2FHxDF%2BO5021pAnSA%3D%3D&sr=1-2-8b2f235a-dddf-4202-bbb9-592393927392"
$ws.Hyperlinks.Add($ws.Range("C15"), $basswoodUrl, "", "", ($basswoodUrl + " "))

# ---------------------------------------------------------------------------
# 3. Update the selected cell shown in the sheet view.
# ---------------------------------------------------------------------------
$ws.Range("I21").Select()
